$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new date headers (BS1, BT1) as literal text shared strings,
# matching how the existing BM1:BR1 headers are stored (t="s", no explicit style).
# Assigning a date-look-alike string straight to .Value auto-converts it to a
# date serial, so instead we push it through as a text formula result and
# flatten it back to a literal value via Copy + PasteSpecial(xlPasteValues).
$ws.Range("BS1").Formula = "=""05/19/2020"""
$ws.Range("BT1").Formula = "=""05/20/2020"""
$ws.Range("BS1:BT1").Copy()
$ws.Range("BS1:BT1").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Fill in the two new cumulative-count columns (BS = 05/19/2020, BT = 05/20/2020)
# for every county row, mirroring the existing BQ/BR pattern.
$ws.Range("BS2").Value = 309
$ws.Range("BT2").Value = 339
$ws.Range("BS3").Value = 7
$ws.Range("BT3").Value = 7
$ws.Range("BS4").Value = 8
$ws.Range("BT4").Value = 8
$ws.Range("BS5").Value = 0
$ws.Range("BT5").Value = 0
$ws.Range("BS6").Value = 0
$ws.Range("BT6").Value = 0
$ws.Range("BS7").Value = 2
$ws.Range("BT7").Value = 2
$ws.Range("BS8").Value = 2
$ws.Range("BT8").Value = 2
$ws.Range("BS9").Value = 15
$ws.Range("BT9").Value = 15
$ws.Range("BS10").Value = 14
$ws.Range("BT10").Value = 14
$ws.Range("BS11").Value = 0
$ws.Range("BT11").Value = 0
$ws.Range("BS12").Value = 1
$ws.Range("BT12").Value = 1
$ws.Range("BS13").Value = 1
$ws.Range("BT13").Value = 2
$ws.Range("BS14").Value = 0
$ws.Range("BT14").Value = 0
$ws.Range("BS15").Value = 0
$ws.Range("BT15").Value = 0
$ws.Range("BS16").Value = 5
$ws.Range("BT16").Value = 5
$ws.Range("BS17").Value = 3
$ws.Range("BT17").Value = 3
$ws.Range("BS18").Value = 0
$ws.Range("BT18").Value = 1
$ws.Range("BS19").Value = 2
$ws.Range("BT19").Value = 2
$ws.Range("BS20").Value = 2
$ws.Range("BT20").Value = 2
$ws.Range("BS21").Value = 0
$ws.Range("BT21").Value = 1
$ws.Range("BS22").Value = 2
$ws.Range("BT22").Value = 2
$ws.Range("BS23").Value = 13
$ws.Range("BT23").Value = 13
$ws.Range("BS24").Value = 0
$ws.Range("BT24").Value = 0
$ws.Range("BS25").Value = 0
$ws.Range("BT25").Value = 0
$ws.Range("BS26").Value = 0
$ws.Range("BT26").Value = 0
$ws.Range("BS27").Value = 0
$ws.Range("BT27").Value = 0
$ws.Range("BS28").Value = 3
$ws.Range("BT28").Value = 4
$ws.Range("BS29").Value = 0
$ws.Range("BT29").Value = 0
$ws.Range("BS30").Value = 0
$ws.Range("BT30").Value = 0
$ws.Range("BS31").Value = 0
$ws.Range("BT31").Value = 0
$ws.Range("BS32").Value = 2
$ws.Range("BT32").Value = 2
$ws.Range("BS33").Value = 3
$ws.Range("BT33").Value = 4
$ws.Range("BS34").Value = 0
$ws.Range("BT34").Value = 0
$ws.Range("BS35").Value = 36
$ws.Range("BT35").Value = 39
$ws.Range("BS36").Value = 0
$ws.Range("BT36").Value = 0
$ws.Range("BS37").Value = 1
$ws.Range("BT37").Value = 1
$ws.Range("BS38").Value = 2
$ws.Range("BT38").Value = 2
$ws.Range("BS39").Value = 1
$ws.Range("BT39").Value = 1
$ws.Range("BS40").Value = 2
$ws.Range("BT40").Value = 2
$ws.Range("BS41").Value = 2
$ws.Range("BT41").Value = 3
$ws.Range("BS42").Value = 2
$ws.Range("BT42").Value = 2
$ws.Range("BS43").Value = 0
$ws.Range("BT43").Value = 0
$ws.Range("BS44").Value = 2
$ws.Range("BT44").Value = 2
$ws.Range("BS45").Value = 10
$ws.Range("BT45").Value = 10
$ws.Range("BS46").Value = 1
$ws.Range("BT46").Value = 2
$ws.Range("BS47").Value = 1
$ws.Range("BT47").Value = 1
$ws.Range("BS48").Value = 422
$ws.Range("BT48").Value = 448

# Mirror the author's final selection/scroll state from the commit.
$ws.Range("BU53").Select()
